# Participation-Tracker.xlsx edit script
# Reproduces the "removing grades file from git" commit:
#  - fills in previously-blank Wheel-3 ("E") participation marks for most students
#  - fixes the F-column participation formula to use a per-row COUNTA(Cn:En)
#    instead of the fixed header range COUNTA($C$1:$E$1)
#  - marks several additional lab/quiz "Y" cells (columns N/O/P)
#  - updates several quiz/exam point totals in the AJ/AK gradebook block
#  - updates the saved selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column E ("Wheel 3") attendance marks that were previously blank
# ---------------------------------------------------------------------
$eValues = @{
    2  = "N"; 3  = "Y"; 4  = "Y"; 5  = "Y"; 6  = "N"; 7  = "Y"; 8  = "Y"; 9  = "Y";
    11 = "Y"; 12 = "Y"; 14 = "N"; 16 = "Y"; 17 = "Y"; 18 = "Y"; 21 = "Y";
    23 = "Y"; 24 = "y"; 25 = "Y"; 26 = "Y"; 30 = "Y"; 31 = "N"; 32 = "Y";
    34 = "Y"; 35 = "Y"; 36 = "N"; 37 = "N"; 39 = "N"; 40 = "Y"; 41 = "N";
    42 = "N"; 43 = "Y"; 44 = "Y"; 45 = "Y"; 46 = "Y";
}
foreach ($r in $eValues.Keys) {
    $ws.Range("E$r").Value = $eValues[$r]
}

# ---------------------------------------------------------------------
# 2. Column F participation-rate formula: use the row's own C:E range for
#    COUNTA instead of the fixed header row reference.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 46; $r++) {
    $ws.Range("F$r").Formula = "=ROUND(COUNTIF(C$r`:E$r, ""Y"")/COUNTA(C$r`:E$r)*5,2)"
}

# ---------------------------------------------------------------------
# 3. Additional "Y" marks in the lab-attendance grid (columns N/O/P)
# ---------------------------------------------------------------------
$opValues = @{
    7  = @("O","P");
    11 = @("O");
    12 = @("O");
    15 = @("O","P");
    16 = @("O","P");
    17 = @("P");
    19 = @("O","P");
    23 = @("O","P");
    25 = @("N","O","P");
    26 = @("O");
    27 = @("O","P");
    30 = @("O","P");
    31 = @("O","P");
    43 = @("O","P");
    44 = @("O");
    45 = @("O","P");
}
foreach ($r in $opValues.Keys) {
    foreach ($col in $opValues[$r]) {
        $ws.Range("$col$r").Value = "Y"
    }
}

# ---------------------------------------------------------------------
# 4. Quiz / exam gradebook updates (columns AI:AT)
# ---------------------------------------------------------------------

# Cells whose "pending entry" yellow highlight is cleared once a final
# score is recorded: copy the plain (unhighlighted) format from AJ28
# before writing the new values.
$ws.Range("AJ28").Copy()
$ws.Range("AJ19").PasteSpecial(-4122)
$ws.Range("AJ20").PasteSpecial(-4122)
$ws.Range("AJ38").PasteSpecial(-4122)
$ws.Range("AJ31").PasteSpecial(-4122)

$ws.Range("AJ16").Value = 14.8
$ws.Range("AJ17").Value = 12.85
$ws.Range("AJ18").Value = 13.95
$ws.Range("AJ19").Value = 15
$ws.Range("AJ20").Value = 12.58
$ws.Range("AJ21").Value = 14.89
$ws.Range("AJ22").Value = 14.95

$ws.Range("AJ28").Value = 8
$ws.Range("AJ29").Value = 7
$ws.Range("AK29").Value = 10
$ws.Range("AL29").Formula = "=AJ29/AK29"
$ws.Range("AJ30").Value = 9
$ws.Range("AJ31").ClearContents()
$ws.Range("AK31").ClearContents()
$ws.Range("AL31").ClearContents()
$ws.Range("AJ32").Value = 9

$ws.Range("AJ37").Value = 58
$ws.Range("AJ38").Value = 50

$ws.Range("AS17").Value = 10

# ---------------------------------------------------------------------
# 5. Restore the saved selection / scroll position
# ---------------------------------------------------------------------
$ws.Range("B19").Select()
$ws.Range("Q32").Select()
